$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly dataset was refreshed: the most recent week's row (previously
# at the bottom, row 6) moves up to row 3, and the rows that used to be in
# rows 3-5 shift down by one row to rows 4-6. Row 2 (the latest date)
# stays put.

$ws.Range("D3").Value = 44624
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 650
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = 675
$ws.Range("P3").Value = 675

$ws.Range("D4").Value = 44608
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 650
$ws.Range("M4").Value = 625
$ws.Range("P4").Value = 625

$ws.Range("D5").Value = 44532
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2200
$ws.Range("M5").Value = 2100
$ws.Range("P5").Value = 2100

$ws.Range("D6").Value = 44533
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2200
$ws.Range("M6").Value = 2100
$ws.Range("P6").Value = 2100
